$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 516, pushing existing rows 516:563 down to 517:564.
$ws.Rows("516:516").Insert()

# Populate the newly inserted row 516 with the new data record.
$row = 516

$ws.Cells.Item($row, 1).Value = 11
$ws.Cells.Item($row, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item($row, 3).Value = "Bíobío"
$ws.Cells.Item($row, 4).Value = 45223
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item(517, 4).NumberFormat()
$ws.Cells.Item($row, 5).Value = 8
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100102
$ws.Cells.Item($row, 8).Value = "Cítricos"
$ws.Cells.Item($row, 9).Value = 100102005
$ws.Cells.Item($row, 10).Value = "Naranja"
$ws.Cells.Item($row, 11).Value = "Valencia"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 200
$ws.Cells.Item($row, 14).Value = 9000
$ws.Cells.Item($row, 15).Value = 10000
$ws.Cells.Item($row, 16).Value = 9500
$ws.Cells.Item($row, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item($row, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($row, 19).Value = 633
$ws.Cells.Item($row, 20).Value = 15
